# repDiffReport.xlsx — re-save/tidy pass picked up by this commit
# ("Add updated fatigueRepReport and MATLAB ppt").
#
# The canonical diff for this workbook is dominated by attributes that
# Excel itself stamps on save (fileVersion/rupBuild, workbookPr
# defaultThemeVersion, the x15ac:absPath scratch value, the x16r2
# mc:Ignorable namespace, sheetView zoomScaleNormal, pageSetup DPI, and
# the bestFit flag / range-coalescing on <cols>). None of those are
# reachable through the Excel object model — they're artifacts of
# whichever build/host wrote the file — so this script drives the part
# of the diff that *is* a real, user-visible action: the author zoomed
# the sheet in, moved the selection, and the column widths for C:L
# collapsed onto a single uniform width.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns C through L all end up at the same rendered width (5 chars).
# Driving ColumnWidth = 4.16666... yields the stored <col width="5".../>
# (Excel's stored width = ColumnWidth + ~0.8333 padding).
$ws.Columns("C:L").ColumnWidth = 4.1666666666666667

# Selection moved from E10 to N12.
$ws.Range("N12").Select()

# Sheet was zoomed way in (205%).
$excel.ActiveWindow.Zoom = 205
